$d = $word.ActiveDocument

# The header table (Tables(1)) in this code-review checklist holds, among
# other fields, the Sprint No. (row 2, col 4) and the Review Date
# (row 3, col 2). Update both to reflect the new sprint / review date.

# Sprint No.: "1" -> "2"
$cellSprint = $d.Tables(1).Cell(2, 4)
$cellSprint.Range.Find.Execute("1", $true, $false, $false, $false, $false, `
    $true, 0, $false, "2", 1)

# Review Date: "02/09/18" -> "02/21/18"
$cellDate = $d.Tables(1).Cell(3, 2)
$cellDate.Range.Find.Execute("02/09/18", $true, $false, $false, $false, $false, `
    $true, 0, $false, "02/21/18", 1)
